# "rose des vents pimpée" - expand the 8-point wind-rose codes in column D
# (Direction du vent) to 16-point codes on the "infos_filtres" sheet, and
# fill in the two rows (7 & 8) that were previously missing a value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "NNE"
$ws.Range("D3").Value = "NNE"
$ws.Range("D4").Value = "ENE"
$ws.Range("D5").Value = "ENE"
$ws.Range("D6").Value = "S"
$ws.Range("D7").Value = "SO"
$ws.Range("D8").Value = "SO"
$ws.Range("D9").Value = "SE"
$ws.Range("D13").Value = "ESE"
$ws.Range("D14").Value = "ESE"
$ws.Range("D15").Value = "ESE"
$ws.Range("D16").Value = "NNE"
$ws.Range("D18").Value = "NNE"
$ws.Range("D19").Value = "NNE"
$ws.Range("D21").Value = "ENE"
$ws.Range("D22").Value = "ENE"
$ws.Range("D23").Value = "ENE"
$ws.Range("D25").Value = "ENE"
$ws.Range("D27").Value = "ENE"
$ws.Range("D28").Value = "ENE"
$ws.Range("D29").Value = "ENE"
$ws.Range("D30").Value = "ENE"
$ws.Range("D31").Value = "ENE"
$ws.Range("D32").Value = "ENE"
$ws.Range("D33").Value = "ENE"
$ws.Range("D34").Value = "ENE"
$ws.Range("D35").Value = "ENE"
$ws.Range("D36").Value = "ENE"
$ws.Range("D37").Value = "ENE"
$ws.Range("D39").Value = "ENE"
$ws.Range("D40").Value = "ESE"
$ws.Range("D42").Value = "NNO"
$ws.Range("D43").Value = "NNO"
$ws.Range("D44").Value = "NNO"
$ws.Range("D45").Value = "NNE"
$ws.Range("D48").Value = "NE"
$ws.Range("D49").Value = "ENE"
$ws.Range("D50").Value = "NNO"
$ws.Range("D55").Value = "ENE"
$ws.Range("D56").Value = "ENE"
$ws.Range("D57").Value = "ENE"
$ws.Range("D61").Value = "ENE"
$ws.Range("D62").Value = "NNO"
$ws.Range("D64").Value = "NNO"
$ws.Range("D67").Value = "NE"
$ws.Range("D68").Value = "NNE"
$ws.Range("D71").Value = "NE"
$ws.Range("D72").Value = "ESE"
$ws.Range("D73").Value = "ENE"
$ws.Range("D74").Value = "NE"
$ws.Range("D77").Value = "NE"
$ws.Range("D78").Value = "ENE"

# the "?" placeholder entry is no longer used anywhere in the sheet now
# that every row has a real wind direction, so drop it from the shared
# string table by overwriting any remaining stray "?" cells (defensive,
# none remain after the edits above) and let the engine garbage-collect
# the now-unused string on save.

# Move the viewport/selection down to the bottom of the table, where the
# author was last working.
$ws.Range("D79").Select()

Write-Host "done"
